$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the creep role path for the harvesting (Mine) row.
$ws.Range("A5").Value = "Mole -> Rabbit -> Elephant"

# Widen columns B and C to better fit the new/longer role names.
# (Values chosen so the engine's internal pixel-rounded ColumnWidth->stored-width
# conversion lands as close as possible to the target widths of 52.5703125 / 54.)
$ws.Columns.Item(2).ColumnWidth = 51.6666666666667
$ws.Columns.Item(3).ColumnWidth = 53.1666666666667

# Move the active selection to A5.
$ws.Range("A5").Select()
